# Minnesota Timberwolves roster: Jaden McDaniels (No. 3) now appears
# ahead of Jaylen Nowell (No. 4) in the roster table — the two players
# swapped roster-table rows (row 3 <-> row 4). Column A (the original
# 0-based row index) keeps its value; every other column (B:K) for the
# two rows is swapped.
#
# Use Range.Copy(destination) rather than reading/writing .Value so that
# each cell's original data type (numbers stay numbers, text stays text
# — e.g. the "Exp" column values "2"/"3" are stored as text, not
# numbers) and style survive the round-trip untouched. A scratch row
# well below the used range holds row 3's data while row 4 is copied
# into row 3; the scratch row is then copied into row 4 and removed so
# the sheet's dimensions/layout end up exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 20

$ws.Range("B3:K3").Copy($ws.Range("B$scratchRow`:K$scratchRow"))
$ws.Range("B4:K4").Copy($ws.Range("B3:K3"))
$ws.Range("B$scratchRow`:K$scratchRow").Copy($ws.Range("B4:K4"))

$ws.Rows("$scratchRow`:$scratchRow").Delete()
